# Insert a new bulleted ("ListParagraph" / numId 2) item right after the
# "My LWJGL debugger for jbullet will probably need to be disabled I
# suspect" bullet, containing:
#   "Indexes for geometries need to have an IntBuffer transport to
#    pipeline! They get copied there any way"
# with "IntBuffer" wrapped in spellStart/spellEnd proofErr markers, same
# as the rest of the document's camel-case / jargon words.

$d = $word.ActiveDocument

# Locate the anchor paragraph and its 1-based index in the Paragraphs
# collection (indexed access gives us a Range with usable Start/End,
# unlike the live enumerator's .Next which this host doesn't resolve).
$anchor = $null
$anchorIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*will probably need to be disabled I suspect*") {
        $anchor = $p
        $anchorIndex = $i
    }
}

if ($anchorIndex -eq -1) {
    throw "Anchor paragraph not found"
}

# Split off a new paragraph after the anchor. Word clones the anchor's
# pPr (ListParagraph style + numPr ilvl=0/numId=2) onto the new, empty
# paragraph.
$anchor.Range.InsertParagraphAfter() | Out-Null

# Re-fetch the freshly created (now next) paragraph by index and grab
# its full range (start through, and including, its paragraph mark) so
# InsertXML replaces the empty run it's holding instead of leaving it
# behind alongside the new content.
$newPara = $d.Paragraphs.Item($anchorIndex + 1)
$newRange = $d.Range($newPara.Range.Start, $newPara.Range.End)

# Use the WordprocessingML "Open XML in Word" package form via InsertXML
# so we control the run boundaries exactly, including the spellStart/
# spellEnd proofErr pair around "IntBuffer" (plain text insertion does
# not synthesize those). The paragraph's own pPr is included so the
# ListParagraph / numbering survive the replace.
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Indexes for geometries need to have an </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>IntBuffer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> transport to pipeline! They get copied there any way</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$newRange.InsertXML($xml) | Out-Null
